$wb = $excel.ActiveWorkbook

# Duplicate the "Turkey" sheet (closest template: same layout/styles) and
# place the copy at the very end of the tab strip, then rename it.
$turkey = $wb.Worksheets.Item("Turkey")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
[void]$turkey.Copy($null, $lastSheet)

$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Update the two cells that differ from the Turkey template.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2474"

# Turkey is no longer the active tab; its selection becomes a full-sheet
# selection (as if the user pressed Ctrl+A there before switching away).
[void]$turkey.Activate()
[void]$turkey.Cells.Select()

# Make the new sheet the active one, with B4 selected (mirrors the
# "last edited cell" selection state captured in the template).
[void]$croatia.Activate()
[void]$croatia.Range("B4").Select()
